# "Generate Report for Handoff" — b.md has been handed off again (new
# source revision 63290e5768f688058c7b37413b0a5c26c308f864) for both the
# zh-cn and de-de targets. Update the Overview rollup plus each language
# sheet's status / latest-handoff-file / latest-handoff-datetime columns
# for the b.md row, and keep the existing hyperlink on the "Latest Handoff
# File" cell pointing at the same link target while refreshing its visible
# text so it mirrors the new cell value.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: row 3 is b.md (row 2 is a.md).
#   B3 = zh-cn status, C3 = de-de status, D3 = Latest Handoff Date
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-24-12 10:24:21"

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 is b.md.
#   C3 = Status
#   D3 = Latest Handoff File (hyperlinked)
#   E3 = Latest Handoff Datetime
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("E3").Value = "2016-03-12 10:24:18"

foreach ($link in $zhcn.Hyperlinks) {
    if ($link.Range.Address() -eq '$D$3') {
        $link.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
    }
}

# ---------------------------------------------------------------------
# de-de sheet: row 3 is b.md.
#   C3 = Status
#   D3 = Latest Handoff File (hyperlinked)
#   E3 = Latest Handoff Datetime
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("E3").Value = "2016-03-12 10:24:21"

foreach ($link in $dede.Hyperlinks) {
    if ($link.Range.Address() -eq '$D$3') {
        $link.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
    }
}
